{"js": "// Update the date line and all 100 two-digit multiplication problems in the\n// table to the new values. Replacements are applied positionally: the Nth\n// paragraph in the document body (date paragraph first, then each table\n// cell's paragraph in row-major order) is matched against the Nth entry of\n// `replacements` below -- this mirrors the structure of the source document\n// (one title paragraph followed by a 20x5 table with a single run per\n// cell) and avoids any ambiguity from duplicate problem text (e.g. \"73\u00d743=\").\nconst replacements = [[\"2023-07-16 Sunday\", \"2023-07-17 Monday\"], [\"53\u00d783=\", \"75\u00d789=\"], [\"76\u00d715=\", \"14\u00d725=\"], [\"10\u00d793=\", \"85\u00d773=\"], [\"93\u00d7100=\", \"42\u00d763=\"], [\"53\u00d770=\", \"48\u00d741=\"], [\"50\u00d729=\", \"54\u00d749=\"], [\"73\u00d743=\", \"60\u00d756=\"], [\"53\u00d730=\", \"67\u00d780=\"], [\"76\u00d710=\", \"83\u00d735=\"], [\"45\u00d731=\", \"41\u00d794=\"], [\"73\u00d743=\", \"58\u00d780=\"], [\"34\u00d734=\", \"15\u00d784=\"], [\"29\u00d737=\", \"27\u00d759=\"], [\"42\u00d721=\", \"28\u00d772=\"], [\"29\u00d786=\", \"37\u00d784=\"], [\"72\u00d729=\", \"35\u00d741=\"], [\"60\u00d794=\", \"33\u00d744=\"], [\"82\u00d768=\", \"43\u00d736=\"], [\"37\u00d723=\", \"42\u00d763=\"], [\"70\u00d739=\", \"13\u00d751=\"], [\"45\u00d725=\", \"75\u00d793=\"], [\"11\u00d760=\", \"35\u00d750=\"], [\"18\u00d736=\", \"46\u00d794=\"], [\"98\u00d782=\", \"31\u00d740=\"], [\"41\u00d770=\", \"40\u00d721=\"], [\"83\u00d7100=\", \"18\u00d794=\"], [\"86\u00d729=\", \"82\u00d758=\"], [\"97\u00d752=\", \"77\u00d765=\"], [\"38\u00d739=\", \"30\u00d730=\"], [\"72\u00d727=\", \"78\u00d712=\"], [\"61\u00d730=\", \"37\u00d764=\"], [\"74\u00d785=\", \"11\u00d764=\"], [\"90\u00d711=\", \"84\u00d729=\"], [\"59\u00d727=\", \"37\u00d796=\"], [\"11\u00d767=\", \"27\u00d719=\"], [\"75\u00d7100=\", \"57\u00d753=\"], [\"84\u00d793=\", \"47\u00d768=\"], [\"59\u00d712=\", \"20\u00d797=\"], [\"54\u00d719=\", \"72\u00d728=\"], [\"50\u00d757=\", \"99\u00d735=\"], [\"92\u00d748=\", \"44\u00d718=\"], [\"90\u00d784=\", \"62\u00d758=\"], [\"22\u00d749=\", \"58\u00d776=\"], [\"93\u00d725=\", \"82\u00d781=\"], [\"40\u00d772=\", \"47\u00d751=\"], [\"62\u00d786=\", \"70\u00d731=\"], [\"65\u00d729=\", \"58\u00d775=\"], [\"71\u00d757=\", \"98\u00d783=\"], [\"56\u00d748=\", \"96\u00d737=\"], [\"73\u00d793=\", \"18\u00d793=\"], [\"10\u00d797=\", \"24\u00d767=\"], [\"85\u00d720=\", \"58\u00d740=\"], [\"58\u00d757=\", \"88\u00d786=\"], [\"73\u00d777=\", \"11\u00d785=\"], [\"22\u00d782=\", \"53\u00d793=\"], [\"17\u00d725=\", \"15\u00d715=\"], [\"75\u00d734=\", \"28\u00d775=\"], [\"22\u00d714=\", \"27\u00d715=\"], [\"26\u00d775=\", \"35\u00d717=\"], [\"83\u00d747=\", \"15\u00d766=\"], [\"60\u00d793=\", \"56\u00d730=\"], [\"77\u00d724=\", \"99\u00d793=\"], [\"44\u00d758=\", \"18\u00d719=\"], [\"68\u00d729=\", \"35\u00d757=\"], [\"62\u00d722=\", \"63\u00d742=\"], [\"91\u00d780=\", \"48\u00d785=\"], [\"12\u00d739=\", \"97\u00d716=\"], [\"86\u00d771=\", \"53\u00d724=\"], [\"45\u00d717=\", \"64\u00d743=\"], [\"16\u00d769=\", \"32\u00d739=\"], [\"41\u00d752=\", \"80\u00d753=\"], [\"85\u00d735=\", \"72\u00d743=\"], [\"28\u00d785=\", \"95\u00d795=\"], [\"23\u00d797=\", \"86\u00d757=\"], [\"81\u00d731=\", \"58\u00d722=\"], [\"11\u00d714=\", \"55\u00d726=\"], [\"63\u00d745=\", \"50\u00d763=\"], [\"35\u00d724=\", \"42\u00d792=\"], [\"53\u00d762=\", \"94\u00d711=\"], [\"97\u00d781=\", \"37\u00d711=\"], [\"59\u00d757=\", \"23\u00d720=\"], [\"96\u00d741=\", \"50\u00d745=\"], [\"97\u00d753=\", \"14\u00d746=\"], [\"76\u00d790=\", \"43\u00d742=\"], [\"77\u00d749=\", \"64\u00d716=\"], [\"62\u00d716=\", \"59\u00d718=\"], [\"13\u00d727=\", \"34\u00d763=\"], [\"89\u00d723=\", \"70\u00d789=\"], [\"42\u00d782=\", \"42\u00d716=\"], [\"73\u00d725=\", \"54\u00d784=\"], [\"100\u00d7100=\", \"36\u00d792=\"], [\"33\u00d742=\", \"69\u00d717=\"], [\"36\u00d796=\", \"14\u00d711=\"], [\"50\u00d783=\", \"86\u00d737=\"], [\"63\u00d733=\", \"15\u00d784=\"], [\"10\u00d747=\", \"28\u00d712=\"], [\"54\u00d787=\", \"46\u00d745=\"], [\"68\u00d755=\", \"54\u00d785=\"], [\"94\u00d733=\", \"52\u00d774=\"], [\"60\u00d772=\", \"14\u00d743=\"]];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} paragraphs, found ${items.length}`\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i}: expected \"${oldText}\" but found \"${para.text}\"`\n    );\n  }\n  if (oldText !== newText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and all 100 two-digit multiplication problems in the\n# table to the new values. The table is 20 rows x 5 columns, one run per\n# cell; cells are addressed by (row, col) so the single duplicated problem\n# text (\"73\u00d743=\") is never ambiguous.\n\n$d = $word.ActiveDocument\n\n$dateOld = \"2023-07-16 Sunday\"\n$dateNew = \"2023-07-17 Monday\"\n$titlePara = $d.Paragraphs(1)\n$titleText = ($titlePara.Range.Text -replace \"[\\r\\x07]\", \"\")\nif ($titleText -ne $dateOld) {\n    throw \"Title paragraph: expected `\"$dateOld`\" but found `\"$titleText`\"\"\n}\n$titlePara.Range.Text = $dateNew\n\n$data = @(\n    @(@(\"53\u00d783=\",\"75\u00d789=\"), @(\"76\u00d715=\",\"14\u00d725=\"), @(\"10\u00d793=\",\"85\u00d773=\"), @(\"93\u00d7100=\",\"42\u00d763=\"), @(\"53\u00d770=\",\"48\u00d741=\")),\n    @(@(\"50\u00d729=\",\"54\u00d749=\"), @(\"73\u00d743=\",\"60\u00d756=\"), @(\"53\u00d730=\",\"67\u00d780=\"), @(\"76\u00d710=\",\"83\u00d735=\"), @(\"45\u00d731=\",\"41\u00d794=\")),\n    @(@(\"73\u00d743=\",\"58\u00d780=\"), @(\"34\u00d734=\",\"15\u00d784=\"), @(\"29\u00d737=\",\"27\u00d759=\"), @(\"42\u00d721=\",\"28\u00d772=\"), @(\"29\u00d786=\",\"37\u00d784=\")),\n    @(@(\"72\u00d729=\",\"35\u00d741=\"), @(\"60\u00d794=\",\"33\u00d744=\"), @(\"82\u00d768=\",\"43\u00d736=\"), @(\"37\u00d723=\",\"42\u00d763=\"), @(\"70\u00d739=\",\"13\u00d751=\")),\n    @(@(\"45\u00d725=\",\"75\u00d793=\"), @(\"11\u00d760=\",\"35\u00d750=\"), @(\"18\u00d736=\",\"46\u00d794=\"), @(\"98\u00d782=\",\"31\u00d740=\"), @(\"41\u00d770=\",\"40\u00d721=\")),\n    @(@(\"83\u00d7100=\",\"18\u00d794=\"), @(\"86\u00d729=\",\"82\u00d758=\"), @(\"97\u00d752=\",\"77\u00d765=\"), @(\"38\u00d739=\",\"30\u00d730=\"), @(\"72\u00d727=\",\"78\u00d712=\")),\n    @(@(\"61\u00d730=\",\"37\u00d764=\"), @(\"74\u00d785=\",\"11\u00d764=\"), @(\"90\u00d711=\",\"84\u00d729=\"), @(\"59\u00d727=\",\"37\u00d796=\"), @(\"11\u00d767=\",\"27\u00d719=\")),\n    @(@(\"75\u00d7100=\",\"57\u00d753=\"), @(\"84\u00d793=\",\"47\u00d768=\"), @(\"59\u00d712=\",\"20\u00d797=\"), @(\"54\u00d719=\",\"72\u00d728=\"), @(\"50\u00d757=\",\"99\u00d735=\")),\n    @(@(\"92\u00d748=\",\"44\u00d718=\"), @(\"90\u00d784=\",\"62\u00d758=\"), @(\"22\u00d749=\",\"58\u00d776=\"), @(\"93\u00d725=\",\"82\u00d781=\"), @(\"40\u00d772=\",\"47\u00d751=\")),\n    @(@(\"62\u00d786=\",\"70\u00d731=\"), @(\"65\u00d729=\",\"58\u00d775=\"), @(\"71\u00d757=\",\"98\u00d783=\"), @(\"56\u00d748=\",\"96\u00d737=\"), @(\"73\u00d793=\",\"18\u00d793=\")),\n    @(@(\"10\u00d797=\",\"24\u00d767=\"), @(\"85\u00d720=\",\"58\u00d740=\"), @(\"58\u00d757=\",\"88\u00d786=\"), @(\"73\u00d777=\",\"11\u00d785=\"), @(\"22\u00d782=\",\"53\u00d793=\")),\n    @(@(\"17\u00d725=\",\"15\u00d715=\"), @(\"75\u00d734=\",\"28\u00d775=\"), @(\"22\u00d714=\",\"27\u00d715=\"), @(\"26\u00d775=\",\"35\u00d717=\"), @(\"83\u00d747=\",\"15\u00d766=\")),\n    @(@(\"60\u00d793=\",\"56\u00d730=\"), @(\"77\u00d724=\",\"99\u00d793=\"), @(\"44\u00d758=\",\"18\u00d719=\"), @(\"68\u00d729=\",\"35\u00d757=\"), @(\"62\u00d722=\",\"63\u00d742=\")),\n    @(@(\"91\u00d780=\",\"48\u00d785=\"), @(\"12\u00d739=\",\"97\u00d716=\"), @(\"86\u00d771=\",\"53\u00d724=\"), @(\"45\u00d717=\",\"64\u00d743=\"), @(\"16\u00d769=\",\"32\u00d739=\")),\n    @(@(\"41\u00d752=\",\"80\u00d753=\"), @(\"85\u00d735=\",\"72\u00d743=\"), @(\"28\u00d785=\",\"95\u00d795=\"), @(\"23\u00d797=\",\"86\u00d757=\"), @(\"81\u00d731=\",\"58\u00d722=\")),\n    @(@(\"11\u00d714=\",\"55\u00d726=\"), @(\"63\u00d745=\",\"50\u00d763=\"), @(\"35\u00d724=\",\"42\u00d792=\"), @(\"53\u00d762=\",\"94\u00d711=\"), @(\"97\u00d781=\",\"37\u00d711=\")),\n    @(@(\"59\u00d757=\",\"23\u00d720=\"), @(\"96\u00d741=\",\"50\u00d745=\"), @(\"97\u00d753=\",\"14\u00d746=\"), @(\"76\u00d790=\",\"43\u00d742=\"), @(\"77\u00d749=\",\"64\u00d716=\")),\n    @(@(\"62\u00d716=\",\"59\u00d718=\"), @(\"13\u00d727=\",\"34\u00d763=\"), @(\"89\u00d723=\",\"70\u00d789=\"), @(\"42\u00d782=\",\"42\u00d716=\"), @(\"73\u00d725=\",\"54\u00d784=\")),\n    @(@(\"100\u00d7100=\",\"36\u00d792=\"), @(\"33\u00d742=\",\"69\u00d717=\"), @(\"36\u00d796=\",\"14\u00d711=\"), @(\"50\u00d783=\",\"86\u00d737=\"), @(\"63\u00d733=\",\"15\u00d784=\")),\n    @(@(\"10\u00d747=\",\"28\u00d712=\"), @(\"54\u00d787=\",\"46\u00d745=\"), @(\"68\u00d755=\",\"54\u00d785=\"), @(\"94\u00d733=\",\"52\u00d774=\"), @(\"60\u00d772=\",\"14\u00d743=\"))\n)\n\n$tbl = $d.Tables(1)\n$r = 1\nforeach ($row in $data) {\n    $c = 1\n    foreach ($pair in $row) {\n        $cell = $tbl.Cell($r, $c)\n        $old = $pair[0]\n        $new = $pair[1]\n        $actual = ($cell.Range.Text -replace \"[\\r\\x07]\", \"\")\n        if ($actual -ne $old) {\n            throw (\"Cell r=\" + $r + \" c=\" + $c + \": expected `\"\" + $old + \"`\" but found `\"\" + $actual + \"`\"\")\n        }\n        if ($old -ne $new) {\n            $cell.Range.Text = $new\n        }\n        $c++\n    }\n    $r++\n}\n\n"}
